$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the meanEMG / legmaxROM values for columns B:E in rows 1-3
$ws.Range("B1").Value = 15
$ws.Range("C1").Value = 16
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 16

$ws.Range("B2").Value = 106.13121722688041
$ws.Range("C2").Value = 104.00990204435887
$ws.Range("D2").Value = 104.92882886823251
$ws.Range("E2").Value = 106.13121722688041

$ws.Range("B3").Value = 105.64292090757452
$ws.Range("C3").Value = 104.46470509868946
$ws.Range("D3").Value = 104.46470509868946
$ws.Range("E3").Value = 105.16444695554922

# Update the active selection to reflect the updated data range
$ws.Range("B1:E3").Select()
